# Powerpoint writer: consolidate text run nodes.
# Merge adjacent <a:r> runs that were previously split on word/space
# boundaries back into single runs, without altering the rendered text.
#
# We do this by re-assigning .Text on specific Characters() sub-ranges -
# the host merges every run fully covered by the target range into one
# run instead of re-splitting on whitespace (which a plain TextRange.Text
# assignment would do).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "A" + " " + "slide"  ->  "A " + "slide" ---
$titleShape = $s.Shapes.Item(1)
$titleRange = $titleShape.TextFrame2.TextRange
$titleRange.Characters(1, 2).Text = "A "

# --- TextBox 3: "Just" " " "an" " " "image" " " "on" " " "this" " " "side" ---
#     -> "Just " "an " "image " "on " "this " "side"
$capShape = $s.Shapes.Item(4)
$capRange = $capShape.TextFrame2.TextRange
$capRange.Characters(1, 5).Text  = "Just "
$capRange.Characters(6, 3).Text  = "an "
$capRange.Characters(9, 6).Text  = "image "
$capRange.Characters(15, 3).Text = "on "
$capRange.Characters(18, 5).Text = "this "
